$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize columns E:G to match the merged column-width group (closest achievable
# to the authored 9.85546875 character width given this host's pixel-grid rounding)
$ws.Columns("E:G").ColumnWidth = 9

# Update cell values (Degradation curves / DC RTE)
$ws.Range("G4").Value = 0.953599999999999
$ws.Range("B5").Value = 0.93514593616264596
$ws.Range("C5").Value = 0.93514593616264596
$ws.Range("D5").Value = 0.950882771660431
$ws.Range("E5").Value = 0.950882771660431
$ws.Range("F5").Value = 0.950882771660431
$ws.Range("G5").Value = 0.95196377677195498
$ws.Range("B6").Value = 0.93373223313476605
$ws.Range("C6").Value = 0.93373223313476605
$ws.Range("D6").Value = 0.94969433141575599
$ws.Range("E6").Value = 0.94969433141575599
$ws.Range("F6").Value = 0.94969433141575599
$ws.Range("G6").Value = 0.950994155345393
$ws.Range("B7").Value = 0.93228887716482101
$ws.Range("C7").Value = 0.93228887716482101
$ws.Range("D7").Value = 0.94874183824396097
$ws.Range("E7").Value = 0.94874183824396097
$ws.Range("F7").Value = 0.94874183824396097
$ws.Range("G7").Value = 0.95019089294109405
$ws.Range("B8").Value = 0.93078427866583802
$ws.Range("C8").Value = 0.93078427866583802
$ws.Range("D8").Value = 0.94793503295484105
$ws.Range("E8").Value = 0.94793503295484105
$ws.Range("F8").Value = 0.94793503295484105
$ws.Range("G8").Value = 0.94950020462213303
$ws.Range("B9").Value = 0.92921336979895497
$ws.Range("C9").Value = 0.92921336979895497
$ws.Range("D9").Value = 0.94724065094697396
$ws.Range("E9").Value = 0.94724065094697396
$ws.Range("F9").Value = 0.94724065094697396
$ws.Range("G9").Value = 0.94888817071072795
$ws.Range("B10").Value = 0.92758288749328999
$ws.Range("C10").Value = 0.92758288749328999
$ws.Range("D10").Value = 0.94663199691967104
$ws.Range("E10").Value = 0.94663199691967104
$ws.Range("F10").Value = 0.94663199691967104
$ws.Range("G10").Value = 0.948336858272536
$ws.Range("B11").Value = 0.92589385433776905
$ws.Range("C11").Value = 0.92589385433776905
$ws.Range("D11").Value = 0.94609368306825403
$ws.Range("E11").Value = 0.94609368306825403
$ws.Range("F11").Value = 0.94609368306825403
$ws.Range("G11").Value = 0.947843862971402
$ws.Range("B12").Value = 0.92415565467384797
$ws.Range("C12").Value = 0.92415565467384797
$ws.Range("D12").Value = 0.94560735245906402
$ws.Range("E12").Value = 0.94560735245906402
$ws.Range("F12").Value = 0.94560735245906402
$ws.Range("G12").Value = 0.94739077026818697
$ws.Range("B13").Value = 0.92238272038232005
$ws.Range("C13").Value = 0.92238272038232005
$ws.Range("D13").Value = 0.94516505311437604
$ws.Range("E13").Value = 0.94516505311437604
$ws.Range("F13").Value = 0.94516505311437604
$ws.Range("G13").Value = 0.94697850254651905
$ws.Range("B14").Value = 0.92058645543171003
$ws.Range("C14").Value = 0.92058645543171003
$ws.Range("D14").Value = 0.944755747752416
$ws.Range("E14").Value = 0.944755747752416
$ws.Range("F14").Value = 0.944755747752416
$ws.Range("G14").Value = 0.94659680971226501
$ws.Range("B15").Value = 0.91876356808079196
$ws.Range("C15").Value = 0.91876356808079196
$ws.Range("D15").Value = 0.94436603844664202
$ws.Range("E15").Value = 0.94436603844664202
$ws.Range("F15").Value = 0.94436603844664202
$ws.Range("G15").Value = 0.94624277082437902
$ws.Range("B16").Value = 0.91694270131199995
$ws.Range("C16").Value = 0.91694270131199995
$ws.Range("D16").Value = 0.94398878474407499
$ws.Range("E16").Value = 0.94398878474407499
$ws.Range("F16").Value = 0.94398878474407499
$ws.Range("G16").Value = 0.94591331998732398
$ws.Range("B17").Value = 0.91510349559824999
$ws.Range("C17").Value = 0.91510349559824999
$ws.Range("D17").Value = 0.94361704167489402
$ws.Range("E17").Value = 0.94361704167489402
$ws.Range("F17").Value = 0.94361704167489402
$ws.Range("G17").Value = 0.94559924726702205
$ws.Range("B18").Value = 0.91328459998148703
$ws.Range("C18").Value = 0.91328459998148703
$ws.Range("D18").Value = 0.94324427329647698
$ws.Range("E18").Value = 0.94324427329647698
$ws.Range("F18").Value = 0.94324427329647698
$ws.Range("G18").Value = 0.94530094826535105
$ws.Range("B19").Value = 0.91147450103442296
$ws.Range("C19").Value = 0.91147450103442296
$ws.Range("D19").Value = 0.94285782645141003
$ws.Range("E19").Value = 0.94285782645141003
$ws.Range("F19").Value = 0.94285782645141003
$ws.Range("G19").Value = 0.94501562908699299
$ws.Range("B20").Value = 0.90969842766601605
$ws.Range("C20").Value = 0.90969842766601605
$ws.Range("D20").Value = 0.94245453834051096
$ws.Range("E20").Value = 0.94245453834051096
$ws.Range("F20").Value = 0.94245453834051096
$ws.Range("G20").Value = 0.94473763747780004
$ws.Range("B21").Value = 0.90794240492960898
$ws.Range("C21").Value = 0.90794240492960898
$ws.Range("D21").Value = 0.94202869950634005
$ws.Range("E21").Value = 0.94202869950634005
$ws.Range("F21").Value = 0.94202869950634005
$ws.Range("G21").Value = 0.94446140514422705
$ws.Range("B22").Value = 0.90624330522344798
$ws.Range("C22").Value = 0.90624330522344798
$ws.Range("D22").Value = 0.94157083965852695
$ws.Range("E22").Value = 0.94157083965852695
$ws.Range("F22").Value = 0.94157083965852695
$ws.Range("G22").Value = 0.94419043568845495
$ws.Range("B23").Value = 0.90458583576141005
$ws.Range("C23").Value = 0.90458583576141005
$ws.Range("D23").Value = 0.94107509372639397
$ws.Range("E23").Value = 0.94107509372639397
$ws.Range("F23").Value = 0.94107509372639397
$ws.Range("G23").Value = 0.94391620689553102
$ws.Range("B24").Value = 0.90298021780622295
$ws.Range("C24").Value = 0.90298021780622295
$ws.Range("D24").Value = 0.940535909879054
$ws.Range("E24").Value = 0.940535909879054
$ws.Range("F24").Value = 0.940535909879054
$ws.Range("G24").Value = 0.943639419103669

# Update active selection to I5 (matches authored sheetView selection)
$ws.Range("I5").Select()
